$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.159.06"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").Value = "3.538.63"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.62"
$ws.Range("E5").Value = "  +6.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.09"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +1.42%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.216"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.57"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("E12").Value = "  -3.46%  "

$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("D14").Value = "4.093.76"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "621.34"
$ws.Range("E15").Value = "  +9.53%  "

$ws.Range("D16").Value = "70.099.50"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.81"
$ws.Range("E17").Value = "  +3.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.97"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("D19").Value = "3.523.15"
$ws.Range("E19").Value = "  -0.41%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  -1.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.50"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.36"
$ws.Range("E23").Value = "  +9.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.70"
$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.00"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("E26").Value = "  +4.12%  "

$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  +8.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.96"
$ws.Range("E29").Value = "  +5.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("E30").Value = "  -2.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.37"
$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.04"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  +18.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.20"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "531.03"
$ws.Range("E36").Value = "  -3.60%  "

$ws.Range("E37").Value = "  -2.33%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.26"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  +7.01%  "

$ws.Range("D41").Value = "0.0₃0779"
$ws.Range("E41").Value = "  -3.23%  "

$ws.Range("D42").Value = "3.532.49"
$ws.Range("E42").Value = "  +3.31%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0465"
$ws.Range("E44").Value = "  +5.38%  "

$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.143"
$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").Value = "  -4.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  -2.90%  "

$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.40"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.82"
$ws.Range("E51").Value = "  -1.24%  "

